# Add carjacking data for 2021-11-29 (extends "through Nov 20" -> "through Nov 21")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet / tab date marker
$ws.Name = "Through 2021-11-21"
$ws.Range("B1").Value = "November 2021 (through November 21)"

# Cell updates: existing values that increment by the new day's counts
$ws.Range("M2").Value  = 14
$ws.Range("X2").Value  = 3
$ws.Range("M3").Value  = 12
$ws.Range("M4").Value  = 9
$ws.Range("B5").Value  = 4
$ws.Range("B6").Value  = 10
$ws.Range("B7").Value  = 3
$ws.Range("M8").Value  = 7
$ws.Range("M9").Value  = 6
$ws.Range("B12").Value = 3
$ws.Range("M13").Value = 3
$ws.Range("B14").Value = 4
$ws.Range("AT14").Value = 1
$ws.Range("M16").Value = 5
$ws.Range("M20").Value = 3
$ws.Range("B21").Value = 1
$ws.Range("M22").Value = 1
$ws.Range("B34").Value = 1
$ws.Range("B35").Value = 2
$ws.Range("BE35").Value = 2
$ws.Range("BE37").Value = 5
$ws.Range("M40").Value = 2
$ws.Range("M42").Value = 1
$ws.Range("X42").Value = 1
$ws.Range("AT52").Value = 1
$ws.Range("AT61").Value = 1
$ws.Range("B64").Value = 5
$ws.Range("B68").Value = 2
$ws.Range("BE69").Value = 1
$ws.Range("B80").Value = 3
$ws.Range("M80").Value = 2
$ws.Range("BE91").Value = 1
